$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value for every data row (rows 2-338).
# All of these cells are being bumped from 46081 to 46082 (one day later).
$range = $ws.Range("C2:C338")
$range.Value = 46082
